# ---------------------------------------------------------------------------
# "removed excess columns in dtr summary and added legends to per employee
#  report"
#
#  1. Clear the stray hard-coded test values that were left in column I
#     (NO. OF OVERTIME HOURS) on rows 7-10 and 17 - they don't belong in a
#     blank template and were removed.
#  2. Add a "Legends:" caption plus three colour-coded legend rows under the
#     DTR summary table explaining the highlight colours used elsewhere in
#     the sheet (request/remark, half-day, absent).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the leftover sample values (1, 1, 1, 1, 0.5) from column I.
#    The cells keep their existing style/border - only the value goes away.
# ---------------------------------------------------------------------------
$ws.Range("I7").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("I17").ClearContents()

# ---------------------------------------------------------------------------
# 2. "Legends:" header - reuse the exact look of the report title
#    (A1 "iRipple, Inc.") by copying its formatting, so no stray new font
#    gets created for this caption.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("E24:P24").Merge()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E24").Value = "Legends:"
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Legend entries - a colour swatch in column E and a bold/underlined
#    description spanning F:P, each swatch/description spanning two rows.
# ---------------------------------------------------------------------------

# --- Legend 1: request / remark (blue, same colour as the OB/offset rows) ---
$ws.Range("E25:E26").Merge()
$ws.Range("F25:P26").Merge()
$ws.Range("E25").Interior.Color = 13411113
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("F25").Font.Bold = $true
$ws.Range("F25").Font.Underline = $true

# --- Legend 2: half-day (orange, same colour as row 17) ---
$ws.Range("E27:E28").Merge()
$ws.Range("F27:P28").Merge()
$ws.Range("E27").Interior.Color = 6737151
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("F27").Font.Bold = $true
$ws.Range("F27").Font.Underline = $true

# --- Legend 3: absent (red, same colour as rows 7-10) ---
$ws.Range("E29:E30").Merge()
$ws.Range("F29:P30").Merge()
$ws.Range("E29").Interior.Color = 6184671
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("F29").Font.Bold = $true
$ws.Range("F29").Font.Underline = $true
